# "Generate Report for Archive"
#
# The underlying localization-status report was regenerated:
#   - the "Status" value for both handback packages flipped from
#     "Ready for handoff" to "In Translation" (Overview sheet columns E/F,
#     and the "Status" column (C) on the per-locale "zh-cn"/"de-de" sheets)
#   - because the status text got shorter, the report generator re-sized
#     (auto-fit) the now-narrower Status columns

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- status text: "Ready for handoff" -> "In Translation" ---------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- narrower Status columns to fit the shorter text ---------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
